$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter -> index map (only the columns we touch)
$colIndex = @{ D = 4; J = 10; K = 11; L = 12; M = 13; P = 16 }

# Target values per row for columns D (date serial), J, K, L, M, P
$data = @{
    2  = @{ D = 44498; J = 40; K = 4000; L = 4000; M = 4000; P = 4000 }
    3  = @{ D = 44504; J = 55 }
    4  = @{ D = 44301; J = 40; K = 3000; L = 3000; M = 3000; P = 3000 }
    5  = @{ D = 44365; K = 5000; L = 5000; M = 5000; P = 5000 }
    6  = @{ D = 44291; J = 35 }
    7  = @{ D = 44312; J = 50 }
    8  = @{ D = 44176; J = 10 }
    9  = @{ D = 44315; J = 40 }
    10 = @{ D = 44313; J = 20 }
    11 = @{ D = 44259; J = 30 }
    12 = @{ D = 44280; J = 55; K = 4000; L = 4000; M = 4000; P = 4000 }
    13 = @{ D = 44497; J = 20 }
    14 = @{ D = 44508; J = 30 }
    15 = @{ D = 44316; J = 20 }
    16 = @{ D = 44509; J = 20 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Cells.Item([int]$row, $colIndex[$col]).Value2 = $cols[$col]
    }
}

$wb.Save()
